$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 82.83048866666665
$ws.Range("H2").Value = 248.491466
$ws.Range("I2").Value = 0.3167437020391103
$ws.Range("J2").Value = 0.3167437020391103
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 1746.216927989611
$ws.Range("R2").Value = 15715.9523519065
$ws.Range("S2").Value = 0.01810526051717649
$ws.Range("T2").Value = 0.01810526051717649

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 82.83048866666665
$ws.Range("H3").Value = 248.491466
$ws.Range("I3").Value = 0.3167437020391103
$ws.Range("J3").Value = 0.3167437020391103
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 24981.68769599931
$ws.Range("R3").Value = 224835.1892639938
$ws.Range("S3").Value = 0.2590170537491782
$ws.Range("T3").Value = 0.2590170537491782

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 82.83048866666665
$ws.Range("H4").Value = 248.491466
$ws.Range("I4").Value = 0.3167437020391103
$ws.Range("J4").Value = 0.3167437020391103
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 3821.40527464867
$ws.Range("R4").Value = 34392.64747183803
$ws.Range("S4").Value = 0.03962138777275566
$ws.Range("T4").Value = 0.03962138777275567

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 152.851481
$ws.Range("H5").Value = 458.554443
$ws.Range("I5").Value = 0.5845039034954311
$ws.Range("J5").Value = 0.5845039034954312
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 3222.38644111604
$ws.Range("R5").Value = 29001.47797004436
$ws.Range("S5").Value = 0.0334105946794316
$ws.Range("T5").Value = 0.0334105946794316

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 152.851481
$ws.Range("H6").Value = 458.554443
$ws.Range("I6").Value = 0.5845039034954311
$ws.Range("J6").Value = 0.5845039034954312
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 46100.02939352018
$ws.Range("R6").Value = 414900.2645416816
$ws.Range("S6").Value = 0.4779778666904217
$ws.Range("T6").Value = 0.4779778666904219

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 152.851481
$ws.Range("H7").Value = 458.554443
$ws.Range("I7").Value = 0.5845039034954311
$ws.Range("J7").Value = 0.5845039034954312
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 7051.841237854757
$ws.Range("R7").Value = 63466.57114069281
$ws.Range("S7").Value = 0.0731154421255778
$ws.Range("T7").Value = 0.07311544212557781

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 25.824378
$ws.Range("H8").Value = 77.47313399999999
$ws.Range("I8").Value = 0.09875239446545848
$ws.Range("J8").Value = 0.0987523944654585
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 544.4247250535659
$ws.Range("R8").Value = 4899.822525482094
$ws.Range("S8").Value = 0.005644746263246413
$ws.Range("T8").Value = 0.005644746263246414

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 25.824378
$ws.Range("H9").Value = 77.47313399999999
$ws.Range("I9").Value = 0.09875239446545848
$ws.Range("J9").Value = 0.0987523944654585
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 7788.636244024195
$ws.Range("R9").Value = 70097.72619621774
$ws.Range("S9").Value = 0.0807547367175792
$ws.Range("T9").Value = 0.08075473671757921

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 25.824378
$ws.Range("H10").Value = 77.47313399999999
$ws.Range("I10").Value = 0.09875239446545848
$ws.Range("J10").Value = 0.0987523944654585
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 1191.414126516374
$ws.Range("R10").Value = 10722.72713864737
$ws.Range("S10").Value = 0.01235291148463288
$ws.Range("T10").Value = 0.01235291148463288
